$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows above the current row 182 (Excel shifts rows
# 182:203 down to 184:205, copying formatting - including the date number
# format on column D - from the row above the insertion point).
$ws.Rows("182:183").Insert()

# New row 182: Kiwi Hayward "Primera"
$ws.Cells.Item(182, 1).Value  = 5
$ws.Cells.Item(182, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(182, 3).Value  = "Maule"
$ws.Cells.Item(182, 4).Value  = 44505
$ws.Cells.Item(182, 5).Value  = 7
$ws.Cells.Item(182, 6).Value  = "Fruta"
$ws.Cells.Item(182, 7).Value  = 100101
$ws.Cells.Item(182, 8).Value  = "Berries"
$ws.Cells.Item(182, 9).Value  = 100101007
$ws.Cells.Item(182, 10).Value = "Kiwi"
$ws.Cells.Item(182, 11).Value = "Hayward"
$ws.Cells.Item(182, 12).Value = "Primera"
$ws.Cells.Item(182, 13).Value = 50
$ws.Cells.Item(182, 14).Value = 15000
$ws.Cells.Item(182, 15).Value = 15000
$ws.Cells.Item(182, 16).Value = 15000
$ws.Cells.Item(182, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(182, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(182, 19).Value = 833
$ws.Cells.Item(182, 20).Value = 18

# New row 183: Kiwi Hayward "Segunda"
$ws.Cells.Item(183, 1).Value  = 5
$ws.Cells.Item(183, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(183, 3).Value  = "Maule"
$ws.Cells.Item(183, 4).Value  = 44505
$ws.Cells.Item(183, 5).Value  = 7
$ws.Cells.Item(183, 6).Value  = "Fruta"
$ws.Cells.Item(183, 7).Value  = 100101
$ws.Cells.Item(183, 8).Value  = "Berries"
$ws.Cells.Item(183, 9).Value  = 100101007
$ws.Cells.Item(183, 10).Value = "Kiwi"
$ws.Cells.Item(183, 11).Value = "Hayward"
$ws.Cells.Item(183, 12).Value = "Segunda"
$ws.Cells.Item(183, 13).Value = 30
$ws.Cells.Item(183, 14).Value = 13000
$ws.Cells.Item(183, 15).Value = 13000
$ws.Cells.Item(183, 16).Value = 13000
$ws.Cells.Item(183, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(183, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(183, 19).Value = 722
$ws.Cells.Item(183, 20).Value = 18
